# Auto-generated edit script: updates FFXIV market-price derived cells
# across multiple worksheets per the scheduled-runner data refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2307.204
$ws.Range("I15").Value = 2307.204
$ws.Range("K15").Value = 6921.612000000001
$ws.Range("M15").Value = -6752.612000000001
$ws.Range("H33").Value = 802.75
$ws.Range("I33").Value = 741.2222
$ws.Range("J33").Value = 853.0909
$ws.Range("K33").Value = 741.2222
$ws.Range("L33").Value = 853.0909
$ws.Range("M33").Value = -512.2222
$ws.Range("N33").Value = -1311.0909
$ws.Range("H95").Value = 34904.168
$ws.Range("J95").Value = 34904.168
$ws.Range("L95").Value = 34904.168
$ws.Range("N95").Value = -40396.168
$ws.Range("H116").Value = 3284.442
$ws.Range("I116").Value = 2807.5925
$ws.Range("J116").Value = 4089.125
$ws.Range("K116").Value = 2807.5925
$ws.Range("L116").Value = 4089.125
$ws.Range("M116").Value = 634.4074999999998
$ws.Range("N116").Value = -10973.125
$ws.Range("H129").Value = 5435809.5
$ws.Range("I129").Value = 41668012
$ws.Range("K129").Value = 125004036
$ws.Range("M129").Value = -124999036
$ws.Range("H132").Value = 3510730.2
$ws.Range("I132").Value = 3705535.5
$ws.Range("K132").Value = 11116606.5
$ws.Range("M132").Value = -11114076.5
$ws.Range("H135").Value = 957.2258
$ws.Range("J135").Value = 1999.6666
$ws.Range("L135").Value = 17996.9994
$ws.Range("N135").Value = -23066.9994
$ws.Range("H141").Value = 351467.25
$ws.Range("I141").Value = 1083.4828
$ws.Range("K141").Value = 3250.4484
$ws.Range("M141").Value = 1929.5516

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 15875
$ws.Range("J24").Value = 15875
$ws.Range("L24").Value = 15875
$ws.Range("N24").Value = -16623
$ws.Range("H34").Value = 27783.666
$ws.Range("I34").Value = 30012.5
$ws.Range("J34").Value = 27146.857
$ws.Range("K34").Value = 30012.5
$ws.Range("L34").Value = 27146.857
$ws.Range("M34").Value = -29741.5
$ws.Range("N34").Value = -27688.857
$ws.Range("H74").Value = 866.1667
$ws.Range("I74").Value = 801.2308
$ws.Range("J74").Value = 1035
$ws.Range("K74").Value = 801.2308
$ws.Range("L74").Value = 1035
$ws.Range("M74").Value = 72.76919999999996
$ws.Range("N74").Value = -2783
$ws.Range("H77").Value = 866.1667
$ws.Range("I77").Value = 801.2308
$ws.Range("J77").Value = 1035
$ws.Range("K77").Value = 4006.154
$ws.Range("L77").Value = 5175
$ws.Range("M77").Value = 361.8459999999995
$ws.Range("N77").Value = -13911
$ws.Range("H100").Value = 15875
$ws.Range("J100").Value = 15875
$ws.Range("L100").Value = 15875
$ws.Range("N100").Value = -18039
$ws.Range("H110").Value = 1574.1428
$ws.Range("I110").Value = 675
$ws.Range("J110").Value = 2773
$ws.Range("K110").Value = 675
$ws.Range("L110").Value = 2773
$ws.Range("M110").Value = 1370
$ws.Range("N110").Value = -6863
$ws.Range("H139").Value = 29539
$ws.Range("J139").Value = 29539
$ws.Range("L139").Value = 29539
$ws.Range("N139").Value = -39819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2161.5833
$ws.Range("I105").Value = 2144.8333
$ws.Range("J105").Value = 2178.3333
$ws.Range("K105").Value = 2144.8333
$ws.Range("L105").Value = 2178.3333
$ws.Range("M105").Value = -397.8332999999998
$ws.Range("N105").Value = -5672.3333
$ws.Range("H134").Value = 2573.8333
$ws.Range("I134").Value = 1798.6666
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 5395.9998
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -2860.9998
$ws.Range("N134").Value = -29070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3786.4883
$ws.Range("I31").Value = 2909.4583
$ws.Range("J31").Value = 4894.316
$ws.Range("K31").Value = 2909.4583
$ws.Range("L31").Value = 4894.316
$ws.Range("M31").Value = -2614.4583
$ws.Range("N31").Value = -5484.316
$ws.Range("H34").Value = 3786.4883
$ws.Range("I34").Value = 2909.4583
$ws.Range("J34").Value = 4894.316
$ws.Range("K34").Value = 2909.4583
$ws.Range("L34").Value = 4894.316
$ws.Range("M34").Value = -2707.4583
$ws.Range("N34").Value = -5298.316
$ws.Range("H58").Value = 11113744
$ws.Range("I58").Value = 1466.7059
$ws.Range("J58").Value = 45460784
$ws.Range("K58").Value = 1466.7059
$ws.Range("L58").Value = 45460784
$ws.Range("M58").Value = -1263.7059
$ws.Range("N58").Value = -45461190
$ws.Range("H127").Value = 32998
$ws.Range("J127").Value = 32998
$ws.Range("L127").Value = 32998
$ws.Range("N127").Value = -42918
$ws.Range("H132").Value = 2985.16
$ws.Range("I132").Value = 2627.7334
$ws.Range("J132").Value = 3521.3
$ws.Range("K132").Value = 7883.2002
$ws.Range("L132").Value = 10563.9
$ws.Range("M132").Value = -5353.2002
$ws.Range("N132").Value = -15623.9
$ws.Range("H136").Value = 11113744
$ws.Range("I136").Value = 1466.7059
$ws.Range("J136").Value = 45460784
$ws.Range("K136").Value = 4400.1177
$ws.Range("L136").Value = 136382352
$ws.Range("M136").Value = -1850.1177
$ws.Range("N136").Value = -136387452
$ws.Range("H137").Value = 35000
$ws.Range("J137").Value = 35000
$ws.Range("L137").Value = 35000
$ws.Range("N137").Value = -45200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 12350
$ws.Range("J87").Value = 15612.5
$ws.Range("L87").Value = 46837.5
$ws.Range("N87").Value = -49333.5
$ws.Range("H90").Value = 12350
$ws.Range("J90").Value = 15612.5
$ws.Range("L90").Value = 140512.5
$ws.Range("N90").Value = -152992.5
$ws.Range("H120").Value = 19800
$ws.Range("I120").Value = 19500
$ws.Range("K120").Value = 58500
$ws.Range("M120").Value = -53662
$ws.Range("H121").Value = 25803.924
$ws.Range("I121").Value = 325
$ws.Range("J121").Value = 30436.455
$ws.Range("K121").Value = 975
$ws.Range("L121").Value = 91309.36500000001
$ws.Range("M121").Value = 335
$ws.Range("N121").Value = -93929.36500000001
$ws.Range("H132").Value = 2535.5715
$ws.Range("I132").Value = 1242.5714
$ws.Range("J132").Value = 3828.5715
$ws.Range("K132").Value = 11183.1426
$ws.Range("L132").Value = 34457.1435
$ws.Range("M132").Value = -8653.142600000001
$ws.Range("N132").Value = -39517.1435
$ws.Range("H133").Value = 7671.6665
$ws.Range("I133").Value = 11343.333
$ws.Range("J133").Value = 4000
$ws.Range("K133").Value = 34029.999
$ws.Range("L133").Value = 12000
$ws.Range("M133").Value = -28969.999
$ws.Range("N133").Value = -22120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 28.52381
$ws.Range("I2").Value = 24.461538
$ws.Range("J2").Value = 35.125
$ws.Range("K2").Value = 24.461538
$ws.Range("L2").Value = 35.125
$ws.Range("M2").Value = 88.538462
$ws.Range("N2").Value = -261.125
$ws.Range("H126").Value = 2649.8845
$ws.Range("I126").Value = 1507.4546
$ws.Range("J126").Value = 3487.6667
$ws.Range("K126").Value = 4522.3638
$ws.Range("L126").Value = 10463.0001
$ws.Range("M126").Value = -2052.3638
$ws.Range("N126").Value = -15403.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4310.857
$ws.Range("I132").Value = 2632.5
$ws.Range("J132").Value = 5836.636
$ws.Range("K132").Value = 7897.5
$ws.Range("L132").Value = 17509.908
$ws.Range("M132").Value = -5367.5
$ws.Range("N132").Value = -22569.908
$ws.Range("H136").Value = 3588.1667
$ws.Range("I136").Value = 2679.7144
$ws.Range("J136").Value = 4860
$ws.Range("K136").Value = 8039.1432
$ws.Range("L136").Value = 14580
$ws.Range("M136").Value = -5489.1432
$ws.Range("N136").Value = -19680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 437290.1
$ws.Range("I122").Value = 668651.6
$ws.Range("J122").Value = 3487.25
$ws.Range("K122").Value = 2005954.8
$ws.Range("L122").Value = 10461.75
$ws.Range("M122").Value = -2003504.8
$ws.Range("N122").Value = -15361.75
$ws.Range("H132").Value = 16409.281
$ws.Range("I132").Value = 3472.5557
$ws.Range("K132").Value = 10417.6671
$ws.Range("M132").Value = -7887.667099999999
